$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

# 1. Replace the "Bouncy ball" example with "Tic Tac Toe"
$ws.Range("C28").Value = "JavaScript Game - Tic Tac Toe"

# 2. Drop the "(Run via Localhost, code in xampp/htdocs/FLATHTML)" note from column D
$ws.Range("D46:D51").ClearContents()

# 3. Append two new example rows, cloning the formatting from the nearest
#    existing rows so the new cells pick up the same styles.

# Row 52 - clone row 51's look (wrapped "Example :" style in column C)
$ws.Range("A51:D51").Copy()
$ws.Range("A52:D52").PasteSpecial(-4122)
$ws.Range("A52").Formula = "=A51+1"
$ws.Range("B52").Value = "php/XML/JavaScript"
$ws.Range("C52").Value = "Example : Merge two XML files into third XML file"
$ws.Range("D52").ClearContents()

# Row 53 - clone row 46's look (plain, non-wrapped column C)
$ws.Range("A46:D46").Copy()
$ws.Range("A53:D53").PasteSpecial(-4122)
$ws.Range("A53").Formula = "=A52+1"
$ws.Range("B53").Value = "JSON/JavaScript"
$ws.Range("C53").Value = "Example : Read External JSON File and display the output"
$ws.Range("D53").ClearContents()

# 4. Leave the view where the author left it
$ws.Activate()
$ws.Range("D39").Select()
